$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Workbook-level changes: rename Sheet1 -> 核心板, add new sheet 底板 after it
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "核心板"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "底板"

# ---------------------------------------------------------------------------
# 核心板: selection moves from C5 (A3:C5) to A1:C15
# ---------------------------------------------------------------------------
$ws1.Range("A1:C15").Select()

# ---------------------------------------------------------------------------
# 底板: build the SD Card / WIFI / LED / ADC-test pin table
# ---------------------------------------------------------------------------

# --- merge the group-label column first, so the subsequent format copy
#     (which carries the exact border/alignment used on 核心板) wins ---
$ws2.Range("A1:A6").Merge()
$ws2.Range("A7:A14").Merge()
$ws2.Range("A15:A20").Merge()
$ws2.Range("A21:A23").Merge()

# --- values ---
$ws2.Range("A1").Value = "SD Card`n需要拉低PC2，不能和WIFI模块共用"
$ws2.Range("B1").Value = "SDIO_D0"
$ws2.Range("C1").Value = "PC8"
$ws2.Range("B2").Value = "SDIO_D1"
$ws2.Range("C2").Value = "PC9"
$ws2.Range("B3").Value = "SDIO_D2"
$ws2.Range("C3").Value = "PC10"
$ws2.Range("B4").Value = "SDIO_D3"
$ws2.Range("C4").Value = "PC11"
$ws2.Range("B5").Value = "SDIO_CLK"
$ws2.Range("C5").Value = "PC12"
$ws2.Range("B6").Value = "SDIO_CMD"
$ws2.Range("C6").Value = "PD2"

$ws2.Range("A7").Value = "WIFI`n不能和SD Card模块共用"
$ws2.Range("B7").Value = "SDIO_D0"
$ws2.Range("C7").Value = "PC8"
$ws2.Range("B8").Value = "SDIO_D1"
$ws2.Range("C8").Value = "PC9"
$ws2.Range("B9").Value = "SDIO_D2"
$ws2.Range("C9").Value = "PC10"
$ws2.Range("B10").Value = "SDIO_D3"
$ws2.Range("C10").Value = "PC11"
$ws2.Range("B11").Value = "SDIO_CLK"
$ws2.Range("C11").Value = "PC12"
$ws2.Range("B12").Value = "SDIO_CMD"
$ws2.Range("C12").Value = "PD2"
$ws2.Range("B13").Value = "WL_REG_ON"
$ws2.Range("C13").Value = "PC2"
$ws2.Range("B14").Value = "WL_HOST_WAKE"
$ws2.Range("C14").Value = "PI11"

$ws2.Range("A21").Value = "LED"
$ws2.Range("B21").Value = "LED_R"
$ws2.Range("C21").Value = "PB0"
$ws2.Range("B22").Value = "LED_G"
$ws2.Range("C22").Value = "PB1"
$ws2.Range("B23").Value = "LED_B"
$ws2.Range("C23").Value = "PA3"

# --- formats, copied from the matching styles already used on 核心板 ---
$ws1.Range("A1").Copy()
$ws2.Range("A1:A5").PasteSpecial(-4122)
$ws2.Range("A7").PasteSpecial(-4122)
$ws2.Range("A15").PasteSpecial(-4122)
$ws2.Range("A21").PasteSpecial(-4122)

$ws1.Range("A2").Copy()
$ws2.Range("A6").PasteSpecial(-4122)
$ws2.Range("A8:A14").PasteSpecial(-4122)
$ws2.Range("A16:A20").PasteSpecial(-4122)
$ws2.Range("A22:A23").PasteSpecial(-4122)

$ws1.Range("B1").Copy()
$ws2.Range("B1:B6").PasteSpecial(-4122)
$ws2.Range("E6").PasteSpecial(-4122)

$ws1.Range("C1").Copy()
$ws2.Range("C1:C6").PasteSpecial(-4122)
$ws2.Range("F6").PasteSpecial(-4122)

$ws1.Range("B3").Copy()
$ws2.Range("B7:B14").PasteSpecial(-4122)
$ws2.Range("B21:B23").PasteSpecial(-4122)

$ws1.Range("C3").Copy()
$ws2.Range("C7:C14").PasteSpecial(-4122)
$ws2.Range("C21:C23").PasteSpecial(-4122)

$ws1.Range("B6").Copy()
$ws2.Range("B15:B20").PasteSpecial(-4122)

$ws1.Range("C6").Copy()
$ws2.Range("C15:C20").PasteSpecial(-4122)

# --- row heights: header rows auto-size back down (no wrap growth),
#     the blank ADC-test rows get the 15.6pt row height used elsewhere ---
$ws2.Rows(1).AutoFit()
$ws2.Rows(7).AutoFit()
$ws2.Range("A15:A20").RowHeight = 15.6

# --- column widths matching the authored sheet ---
$ws2.Columns("A").ColumnWidth = 8.5546875
$ws2.Columns("B").ColumnWidth = 16.33203125
$ws2.Columns("C").ColumnWidth = 10.109375

# --- page setup to match 核心板 ---
$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

# --- selection on the new sheet ---
$ws2.Range("D14").Select()

# keep 核心板 as the active/visible tab, as in the source workbook
$ws1.Activate()
